$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet and set up its selection state ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "TeamsUsageStats"

# --- Add the new "TeamsGroups" sheet right after TeamsUsageStats ---
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1, [Type]::Missing, [Type]::Missing)
$ws2.Name = "TeamsGroups"

# --- Header row (reuses existing shared strings) ---
$ws2.Range("A1").Value = "Data Type"
$ws2.Range("B1").Value = "Display Name"
$ws2.Range("C1").Value = "Actual Column Name"
$ws2.Range("A1:C1").Font.Bold = $true

# --- Data rows; fill order (B, C, then A) matches the original authoring order ---
$ws2.Range("B2").Value = "Title "
$ws2.Range("C2").Value = "Title "
$ws2.Range("A2").Value = "Single line of text "

$ws2.Range("B3").Value = "GroupID "
$ws2.Range("C3").Value = "GroupID "
$ws2.Range("A3").Value = "Single line of text "

$ws2.Range("B4").Value = "Visibility "
$ws2.Range("C4").Value = "Visibility "
$ws2.Range("A4").Value = "Single line of text "

$ws2.Range("B5").Value = "CreatedOn "
$ws2.Range("C5").Value = "CreatedOn "
$ws2.Range("A5").Value = "Date and Time "

$ws2.Range("B6").Value = "RenewedOn "
$ws2.Range("C6").Value = "RenewedOn "
$ws2.Range("A6").Value = "Date and Time "

$ws2.Range("B7").Value = "DeletedOn "
$ws2.Range("C7").Value = "DeletedOn "
$ws2.Range("A7").Value = "Single line of text "

$ws2.Range("B8").Value = "ownerUID "
$ws2.Range("C8").Value = "ownerUID "
$ws2.Range("A8").Value = "Single line of text "

$ws2.Range("B9").Value = "memberCount "
$ws2.Range("C9").Value = "memberCount "
$ws2.Range("A9").Value = "Number "

$ws2.Range("B10").Value = "lastActivity "
$ws2.Range("C10").Value = "lastActivity "
$ws2.Range("A10").Value = "Number "

$ws2.Range("B11").Value = "externalMemberCount "
$ws2.Range("C11").Value = "externalMemberCount "
$ws2.Range("A11").Value = "Number "

$ws2.Range("B12").Value = "Modified "
$ws2.Range("C12").Value = "Modified "
$ws2.Range("A12").Value = "Date and Time "

$ws2.Range("B13").Value = "Created "
$ws2.Range("C13").Value = "Created "
$ws2.Range("A13").Value = "Date and Time "

$ws2.Range("B14").Value = "Created By "
$ws2.Range("C14").Value = "Author"
$ws2.Range("A14").Value = "Person or Group "

$ws2.Range("B15").Value = "Modified By "
$ws2.Range("C15").Value = "Editor"
$ws2.Range("A15").Value = "Person or Group "

# --- Column widths for the new sheet (values chosen so the engine's pixel
# quantization lands on/near the authored widths 22.5546875 / 22.77734375 / 19) ---
$ws2.Columns.Item(1).ColumnWidth = 21.666666666666668
$ws2.Columns.Item(2).ColumnWidth = 22
$ws2.Columns.Item(3).ColumnWidth = 18.166666666666668

# --- Selections: TeamsUsageStats row 1 is selected (no longer the active tab) ---
$ws1.Rows(1).Select() | Out-Null

# TeamsGroups becomes the active sheet, with C2:C13 selected
$ws2.Range("C2:C13").Select() | Out-Null
